# Replace the separate "Protein" / "Ontogeny" columns with a single
# "Protein Ontogenies" column on the IndividualBiometrics sheet:
#   - H1 header text "Protein"  -> "Protein Ontogenies"
#   - I1 header "Ontogeny" column is removed (no longer used)
#   - H2 gets the new value "CYP3A4:CYP3A4,CYP2D6:CYP2C8"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IndividualBiometrics")

# Rename the "Protein" header to the merged "Protein Ontogenies" header.
$ws.Range("H1").Value = "Protein Ontogenies"

# Drop the old "Ontogeny" header/column - it is no longer needed.
$ws.Range("I1").ClearContents()

# Populate the new merged column with the protein ontogeny mapping value.
$ws.Range("H2").Value = "CYP3A4:CYP3A4,CYP2D6:CYP2C8"

# Restore the sheet's active selection.
$ws.Activate()
$ws.Range("Q6").Select() | Out-Null
